# Updated cryptos list on Thu Oct 24 10:59:35 UTC 2024 with GitHub Actions
#
# The "Price" column (D) stores purely decimal-looking numbers (e.g. "588.65")
# as TEXT in the workbook (Excel shows "8.09", "69.60" etc. with the exact
# number of decimals the scraper produced). Assigning such a string straight
# to Range.Value would make Excel auto-coerce it to a real number (losing the
# trailing zero / exact text form), so those cells are written with a leading
# apostrophe, which is the standard COM/UI idiom for "force text".
# Values that already fail to parse as a plain number (e.g. "66.980.42",
# thousands-dotted prices, or the subscript-digit PEPE/BabyDogeCoin prices)
# are left alone since Excel keeps them as text natively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value  = '66.980.42'
$ws.Range('E2').Value  = '  +0.85%  '

$ws.Range('D3').Value  = '2.525.39'
$ws.Range('E3').Value  = '  -2.13%  '

$ws.Range('E4').Value  = '  -0.04%  '

$ws.Range('D5').Value  = "'588.65"
$ws.Range('E5').Value  = '  +1.04%  '

$ws.Range('D6').Value  = "'172.82"
$ws.Range('E6').Value  = '  +4.26%  '

$ws.Range('E7').Value  = '  +0.01%  '

$ws.Range('E8').Value  = '  +0.13%  '

$ws.Range('D9').Value  = '2.525.60'
$ws.Range('E9').Value  = '  -2.07%  '

$ws.Range('D10').Value = "'0.137"
$ws.Range('E10').Value = '  +0.47%  '

$ws.Range('E11').Value = '  +2.13%  '

$ws.Range('E12').Value = '  -0.38%  '

$ws.Range('E13').Value = '  -3.72%  '

$ws.Range('D14').Value = "'26.49"
$ws.Range('E14').Value = '  -0.62%  '

$ws.Range('D15').Value = '2.985.65'
$ws.Range('E15').Value = '  -2.13%  '

$ws.Range('E16').Value = '  -1.08%  '

$ws.Range('D17').Value = '66.835.50'
$ws.Range('E17').Value = '  +0.61%  '

$ws.Range('D18').Value = '2.494.48'
$ws.Range('E18').Value = '  -4.86%  '

$ws.Range('D19').Value = "'8.09"
$ws.Range('E19').Value = '  +5.02%  '

$ws.Range('D20').Value = "'11.29"
$ws.Range('E20').Value = '  -0.96%  '

$ws.Range('D21').Value = "'354.65"
$ws.Range('E21').Value = '  +0.79%  '

$ws.Range('E22').Value = '  -1.10%  '

$ws.Range('E23').Value = '  +0.29%  '

$ws.Range('E24').Value = '  +5.61%  '

$ws.Range('E25').Value = '  +0.02%  '

$ws.Range('D26').Value = "'69.60"
$ws.Range('E26').Value = '  +1.34%  '

$ws.Range('D27').Value = "'9.92"
$ws.Range('E27').Value = '  -0.68%  '

$ws.Range('D28').Value = "'0.998"
$ws.Range('E28').Value = '  -0.28%  '

$ws.Range('D30').Value = '0.0₃0973'
$ws.Range('E30').Value = '  -0.87%  '

$ws.Range('D31').Value = "'530.52"
$ws.Range('E31').Value = '  -0.34%  '

$ws.Range('D32').Value = "'8.11"
$ws.Range('E32').Value = '  +1.26%  '

$ws.Range('E33').Value = '  -0.11%  '

$ws.Range('E35').Value = '  -0.71%  '

$ws.Range('E36').Value = '  +0.00%  '

$ws.Range('E37').Value = '  -0.19%  '

$ws.Range('D38').Value = "'156.84"
$ws.Range('E38').Value = '  +0.21%  '

$ws.Range('D39').Value = "'18.56"
$ws.Range('E39').Value = '  -0.78%  '

$ws.Range('D40').Value = "'18.43"
$ws.Range('E40').Value = '  +1.02%  '

$ws.Range('E41').Value = '  -1.96%  '

$ws.Range('E42').Value = '  +0.46%  '

$ws.Range('E43').Value = '  +0.37%  '

$ws.Range('E44').Value = '  -0.03%  '

$ws.Range('E45').Value = '  +3.28%  '

$ws.Range('D46').Value = "'148.92"
$ws.Range('E46').Value = '  +0.02%  '

$ws.Range('E47').Value = '  -2.00%  '

$ws.Range('D48').Value = '0.0₆0276'

$ws.Range('E49').Value = '  -0.86%  '

$ws.Range('D50').Value = "'1.69"
$ws.Range('E50').Value = '  -1.20%  '

$ws.Range('E51').Value = '  -0.18%  '
